$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.259.72"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "'3.888.47"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'482.42"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").Value = "'145.33"
$ws.Range("E6").Value = "  -1.45%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("D10").Value = "'0.181"
$ws.Range("E10").Value = "  +7.32%  "
$ws.Range("D11").Value = "'0.0000353"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "'43.07"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "'4.504.36"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("D15").Value = "'3.878.01"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "'14.21"
$ws.Range("E16").Value = "  -2.50%  "
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "'19.92"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "'68.279.99"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "'429.42"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("E22").Value = "  +8.24%  "
$ws.Range("D23").Value = "'14.82"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").Value = "'12.45"
$ws.Range("E24").Value = "  +16.59%  "
$ws.Range("D25").Value = "'88.83"
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("D26").Value = "'3.66"
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D27").Value = "'11.01"
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("D28").Value = "'37.28"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D30").Value = "'715.60"
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").Value = "'13.48"
$ws.Range("E31").Value = "  +1.97%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("E33").Value = "  +2.90%  "
$ws.Range("D34").Value = "'61.93"
$ws.Range("E34").Value = "  +5.95%  "
$ws.Range("D35").Value = "'0.0₃0878"
$ws.Range("E35").Value = "  -2.71%  "
$ws.Range("D36").Value = "'6.06"
$ws.Range("E36").Value = "  +9.11%  "
$ws.Range("D37").Value = "'40.90"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("D38").Value = "'0.399"
$ws.Range("E38").Value = "  +16.13%  "
$ws.Range("E39").Value = "  -3.22%  "
$ws.Range("D40").Value = "'3.02"
$ws.Range("E40").Value = "  +6.38%  "
$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "'0.0498"
$ws.Range("E42").Value = "  +6.05%  "
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("E44").Value = "  -3.04%  "
$ws.Range("E45").Value = "  +2.05%  "
$ws.Range("E46").Value = "  +4.21%  "
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D48").Value = "'3.37"
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.0₆0347"
$ws.Range("E49").Value = "  +26.55%  "
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("D51").Value = "'144.37"
$ws.Range("E51").Value = "  -2.40%  "
